# Apply the AutoModesAnalysis.xlsx edits:
#  - Rename Sheet1 -> Summary
#  - Add two new sheets: c_switch_r, c_switch_l (with data tables)
#  - Make c_switch_l the active/selected tab

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "Summary" ---
$summary = $wb.Worksheets.Item(1)
$summary.Name = "Summary"

# --- Add "c_switch_r" sheet (after Summary) ---
$switchR = $wb.Worksheets.Add($null, $summary)
$switchR.Name = "c_switch_r"

$rData = @(
    @(1, "Y", 2.6),
    @(2, "Y", 2.4),
    @(3, "Y", 2.4),
    @(4, "Y", 2.4),
    @(5, "Y", 2.4),
    @(6, "Y", 2.6),
    @(7, "Y", 2.5),
    @(8, "Y", 2.5),
    @(9, "Y", 2.5),
    @(10, "Y", 2.5)
)

$switchR.Cells.Item(1, 1).Value = "Attempt"
$switchR.Cells.Item(1, 2).Value = "Success"

for ($i = 0; $i -lt $rData.Count; $i++) {
    $row = $i + 2
    $switchR.Cells.Item($row, 1).Value = $rData[$i][0]
    $switchR.Cells.Item($row, 2).Value = $rData[$i][1]
}

$switchR.Cells.Item(1, 3).Value = "Time (s)"

for ($i = 0; $i -lt $rData.Count; $i++) {
    $row = $i + 2
    $switchR.Cells.Item($row, 3).Value = $rData[$i][2]
}

$switchR.Range("B2:B11").HorizontalAlignment = -4108
$switchR.Range("A1:C11").Select() | Out-Null

# --- Add "c_switch_l" sheet (after c_switch_r) ---
$switchL = $wb.Worksheets.Add($null, $switchR)
$switchL.Name = "c_switch_l"

$lData = @(
    @(1, "Y", 7.2),
    @(2, "Y", 7.2),
    @(3, "Y", 7.2),
    @(4, "Y", 7.2),
    @(5, "Y", 8.2),
    @(6, "Y", 8.2),
    @(7, "Y", 7.9),
    @(8, "Y", 7.3),
    @(9, "Y", 7.3),
    @(10, "Y", 7.3)
)

$switchL.Cells.Item(1, 1).Value = "Attempt"
$switchL.Cells.Item(1, 2).Value = "Success"
$switchL.Cells.Item(1, 3).Value = "Time (s)"

for ($i = 0; $i -lt $lData.Count; $i++) {
    $row = $i + 2
    $switchL.Cells.Item($row, 1).Value = $lData[$i][0]
    $switchL.Cells.Item($row, 2).Value = $lData[$i][1]
    $switchL.Cells.Item($row, 3).Value = $lData[$i][2]
}

$switchL.Range("B2:B11").HorizontalAlignment = -4108

# --- Select C12 on c_switch_l and make it the active sheet/tab ---
$switchL.Activate() | Out-Null
$switchL.Range("C12").Select() | Out-Null
